$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 75, shifting existing rows 75-131 down to 76-132
$ws.Rows.Item(75).EntireRow.Insert()

# Populate the newly inserted row 75 with the new weekly record
$ws.Range("A75").Value = 7
$ws.Range("B75").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C75").Value = 'Ñuble'
$ws.Range("D75").Value = 45040
$ws.Range("E75").Value = 16
$ws.Range("F75").Value = 100112037
$ws.Range("G75").Value = 'Cebollín'
$ws.Range("H75").Value = 'Sin especificar'
$ws.Range("I75").Value = 'Primera'
$ws.Range("J75").Value = 60
$ws.Range("K75").Value = 7000
$ws.Range("L75").Value = 7000
$ws.Range("M75").Value = 7000
$ws.Range("N75").Value = '$/paquete 36 unidades'
$ws.Range("O75").Value = 'Provincia de Diguillín'
$ws.Range("P75").Value = 194
$ws.Range("Q75").Value = 36
$ws.Range("R75").Value = 'Hortaliza'
